# "fix typo of indexes" - updates recalculated values in the data table
# after correcting an off-by-index bug in how rows referenced their source data.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 0.72155591310289047
$ws.Range("G2").Value = 0.51783640594870639

$ws.Range("E3").Value = 0.73033464779552981
$ws.Range("G3").Value = 0.48226681424374995

$ws.Range("E4").Value = 0.73890647070283499
$ws.Range("G4").Value = 0.4570670737329014

$ws.Range("E5").Value = 0.74171435119128082
$ws.Range("G5").Value = 0.43817092763914656

$ws.Range("E6").Value = 0.7493388853304237
$ws.Range("G6").Value = 0.42330660901653777

$ws.Range("F14").Value = 0.9845849463137617
$ws.Range("H14").Value = 0.74608591486218567

$ws.Range("E15").Value = 0.99829545665405295
$ws.Range("H15").Value = 0.74957313655838465

$ws.Range("F20").Value = 0.95828872109608543
$ws.Range("H20").Value = 0.73911907993526249

$ws.Range("F22").Value = 0.89583167445557765
$ws.Range("H22").Value = 0.72096257493419136
